$data = @(
  @("2024-09-15",9,0.087),
  @("2024-09-15",10,0.08),
  @("2024-09-15",11,0.078),
  @("2024-09-15",12,0.077),
  @("2024-09-15",13,0.078),
  @("2024-09-15",14,0.076),
  @("2024-09-15",15,0.075),
  @("2024-09-15",16,0.073),
  @("2024-09-15",17,0.073),
  @("2024-09-15",18,0.071),
  @("2024-09-15",19,0.073),
  @("2024-09-15",20,0.075),
  @("2024-09-15",21,0.075),
  @("2024-09-15",22,0.072),
  @("2024-09-15",23,0.07),
  @("2024-09-16",0,0.088),
  @("2024-09-16",1,0.089),
  @("2024-09-16",2,0.091),
  @("2024-09-16",3,0.093),
  @("2024-09-16",4,0.098),
  @("2024-09-16",5,0.107),
  @("2024-09-16",6,0.13),
  @("2024-09-16",7,0.177),
  @("2024-09-16",8,0.213),
  @("2024-09-16",9,0.238),
  @("2024-09-16",10,0.239),
  @("2024-09-16",11,0.217),
  @("2024-09-16",12,0.212),
  @("2024-09-16",13,0.205),
  @("2024-09-16",14,0.208),
  @("2024-09-16",15,0.198),
  @("2024-09-16",16,0.194),
  @("2024-09-16",17,0.18),
  @("2024-09-16",18,0.171),
  @("2024-09-16",19,0.172),
  @("2024-09-16",20,0.188),
  @("2024-09-16",21,0.163),
  @("2024-09-16",22,0.14),
  @("2024-09-16",23,0.134),
  @("2024-09-17",0,0.112),
  @("2024-09-17",1,0.109),
  @("2024-09-17",2,0.104),
  @("2024-09-17",3,0.105),
  @("2024-09-17",4,0.106),
  @("2024-09-17",5,0.114),
  @("2024-09-17",6,0.151),
  @("2024-09-17",7,0.192),
  @("2024-09-17",8,0.217),
  @("2024-09-17",9,0.229),
  @("2024-09-17",10,0.235),
  @("2024-09-17",11,0.228),
  @("2024-09-17",12,0.229),
  @("2024-09-17",13,0.226),
  @("2024-09-17",14,0.222),
  @("2024-09-17",15,0.218),
  @("2024-09-17",16,0.216),
  @("2024-09-17",17,0.202),
  @("2024-09-17",18,0.187),
  @("2024-09-17",19,0.183),
  @("2024-09-17",20,0.186),
  @("2024-09-17",21,0.186),
  @("2024-09-17",22,0.163),
  @("2024-09-17",23,0.134),
  @("2024-09-18",0,0.12),
  @("2024-09-18",1,0.118),
  @("2024-09-18",2,0.11),
  @("2024-09-18",3,0.11),
  @("2024-09-18",4,0.111),
  @("2024-09-18",5,0.121),
  @("2024-09-18",6,0.163),
  @("2024-09-18",7,0.204),
  @("2024-09-18",8,0.236),
  @("2024-09-18",9,0.244),
  @("2024-09-18",10,0.247),
  @("2024-09-18",11,0.244),
  @("2024-09-18",12,0.247),
  @("2024-09-18",13,0.264),
  @("2024-09-18",14,0.246),
  @("2024-09-18",15,0.243),
  @("2024-09-18",16,0.235),
  @("2024-09-18",17,0.224),
  @("2024-09-18",18,0.212),
  @("2024-09-18",19,0.204),
  @("2024-09-18",20,0.197),
  @("2024-09-18",21,0.187),
  @("2024-09-18",22,0.165),
  @("2024-09-18",23,0.136),
  @("2024-09-19",0,0.117),
  @("2024-09-19",1,0.113),
  @("2024-09-19",2,0.109),
  @("2024-09-19",3,0.109),
  @("2024-09-19",4,0.11),
  @("2024-09-19",5,0.12),
  @("2024-09-19",6,0.159),
  @("2024-09-19",7,0.201),
  @("2024-09-19",8,0.23),
  @("2024-09-19",9,0.235),
  @("2024-09-19",10,0.24),
  @("2024-09-19",11,0.232),
  @("2024-09-19",12,0.238),
  @("2024-09-19",13,0.248),
  @("2024-09-19",14,0.23),
  @("2024-09-19",15,0.23),
  @("2024-09-19",16,0.222),
  @("2024-09-19",17,0.215),
  @("2024-09-19",18,0.197),
  @("2024-09-19",19,0.185),
  @("2024-09-19",20,0.182),
  @("2024-09-19",21,0.179),
  @("2024-09-19",22,0.159),
  @("2024-09-19",23,0.131),
  @("2024-09-20",0,0.111),
  @("2024-09-20",1,0.107),
  @("2024-09-20",2,0.105),
  @("2024-09-20",3,0.105),
  @("2024-09-20",4,0.105),
  @("2024-09-20",5,0.116),
  @("2024-09-20",6,0.146),
  @("2024-09-20",7,0.193),
  @("2024-09-20",8,0.213),
  @("2024-09-20",9,0.224),
  @("2024-09-20",10,0.228),
  @("2024-09-20",11,0.221),
  @("2024-09-20",12,0.229),
  @("2024-09-20",13,0.226),
  @("2024-09-20",14,0.213),
  @("2024-09-20",15,0.205),
  @("2024-09-20",16,0.193),
  @("2024-09-20",17,0.187),
  @("2024-09-20",18,0.177),
  @("2024-09-20",19,0.167),
  @("2024-09-20",20,0.152),
  @("2024-09-20",21,0.152),
  @("2024-09-20",22,0.129),
  @("2024-09-20",23,0.108),
  @("2024-09-21",0,0.093),
  @("2024-09-21",1,0.092),
  @("2024-09-21",2,0.091),
  @("2024-09-21",3,0.09),
  @("2024-09-21",4,0.091),
  @("2024-09-21",5,0.094),
  @("2024-09-21",6,0.096),
  @("2024-09-21",7,0.099),
  @("2024-09-21",8,0.099),
  @("2024-09-21",9,0.103),
  @("2024-09-21",10,0.106),
  @("2024-09-21",11,0.106),
  @("2024-09-21",12,0.106),
  @("2024-09-21",13,0.123),
  @("2024-09-21",14,0.122),
  @("2024-09-21",15,0.124),
  @("2024-09-21",16,0.116),
  @("2024-09-21",17,0.118),
  @("2024-09-21",18,0.112),
  @("2024-09-21",19,0.097),
  @("2024-09-21",20,0.094),
  @("2024-09-21",21,0.088),
  @("2024-09-21",22,0.085),
  @("2024-09-21",23,0.082),
  @("2024-09-22",0,0.081),
  @("2024-09-22",1,0.08),
  @("2024-09-22",2,0.077),
  @("2024-09-22",3,0.078),
  @("2024-09-22",4,0.077),
  @("2024-09-22",5,0.081),
  @("2024-09-22",6,0.082),
  @("2024-09-22",7,0.083),
  @("2024-09-22",8,0.086),
  @("2024-09-22",9,0.084)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = 2 + $i
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = "dd.mm.yyyy"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
